# Add evaluation for topologies and plots
$wb = $excel.ActiveWorkbook

# --- Rename "Sheet2" (4th tab) to "PCIe info" and populate it ---
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "PCIe info"

# PCIe link numbers
$ws.Range("A2").Value = "Link Speed"
$ws.Range("B2").Value = 985

$ws.Range("A3").Value = "Lanes"
$ws.Range("B3").Value = 8
$ws.Range("I3").Formula = "=40/8"

$ws.Range("A4").Value = "Peak"
$ws.Range("B4").Formula = "=B2*B3"

$ws.Range("A5").Value = " "
$ws.Range("A5").WrapText = $true
$ws.Range("B5").Formula = "=1/((1/7.88)+(1/12.5)+(1/7.88))"

# Peak bandwidth comparison table
$ws.Range("E8").Value = "Configuration"
$ws.Range("F8").Value = "Peak Bandwidth (GB/s)"

$ws.Range("E9").Value = "MPI+PCIe"
$ws.Range("F9").Value = 2.995

$ws.Range("E10").Value = "Within Node 1 Channel"
$ws.Range("F10").Value = 5

$ws.Range("E11").Value = "Within Node 4 Channel"
$ws.Range("F11").Value = 20

$ws.Range("E12").Value = "Fully Connect"
$ws.Range("F12").Value = 15

# Column widths to fit the new content
# (values back-solved so the engine's internal 1/6-char pixel-grid rounding
# lands as close as possible to the target stored widths 12.21875/20/25.77734375)
$ws.Columns.Item(1).ColumnWidth = 11.333333333333332
$ws.Columns.Item(5).ColumnWidth = 19.166666666666664
$ws.Columns.Item(6).ColumnWidth = 25.0

# Make "PCIe info" the active sheet/tab (this also clears tabSelected from
# whichever sheet - Sheet3 - previously held it)
$ws.Activate()

# Selection on the new sheet
$ws.Range("E8:F12").Select()
